$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the tiny floating point precision of the existing A9 timestamp
$ws.Range("A9").Value = 45866.41686357639

# Append the new row of data (row 10)
$ws.Range("A10").Value = 45866.45858155453
$ws.Range("B10").Value = 2025
$ws.Range("C10").Value = 31
$ws.Range("D10").Value = 17.17
$ws.Range("E10").Value = 81.16
$ws.Range("F10").Value = 595.9
$ws.Range("G10").Value = 14.17
$ws.Range("H10").Value = "ESE"
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = "11:00:21"

# Match the style used by the other date cells in column A
$ws.Range("A10").NumberFormat = $ws.Range("A9").NumberFormat
